# Add a new "PEAAC" (Process Emissions Additions and Costs) acronym row to the
# "Key to Variables" sheet, directly above the existing "PERAC" row (old row 188),
# which shifts all subsequent rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Key to Variables")
$ws.Activate()

# Insert a new blank row at position 188 (pushes old row 188 onward down to 189+)
$ws.Rows.Item(188).Insert()

# Populate the new row's cells
$ws.Range("A188").Value = "indst"
$ws.Range("B188").Value = "PEAAC"
$ws.Range("C188").Value = "Process Emissions Additions and Costs"
$ws.Range("D188").Value = "Process Emissions Additions and Costs, Marginal Cost Definitions"
$ws.Range("F188").Value = "optional"
$ws.Range("G188").Value = "You want to test the addition rather than the reduction of process emissions for select industries"

# Match the formatting used by other rows whose "Importance to Update" (column F)
# value is "optional" (style used in e.g. F3), so F188 gets the correct fill/style.
$ws.Range("F3").Copy()
$ws.Range("F188").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# This row wraps onto two lines like its neighbours, so match their row height.
$ws.Rows.Item(188).RowHeight = 30

# Update the visible selection to reflect where the edit was made.
$ws.Range("G189").Select()
